# This edit normalises the "id" column and the primary-key constraint of
# six tables in the create-tables script (make, belong_to, collect,
# comment, open, listen):
#
#   1. "    id bigint,"  ->  "    id bigint not null,"
#   2. each table's (sometimes composite / differently named) primary
#      key is simplified down to "primary key (id)".
#
# The other tables (sys_user, artist, album, track) already read
# "id bigint not null," and are left untouched.

$d = $word.ActiveDocument

# --- 1) id column: mark "not null" -----------------------------------
# Walk every "id bigint," occurrence in the document and turn the
# trailing comma-run into " not null," in place, so the existing
# "bigint" run/proofErr markers are left alone (only six such
# paragraphs exist; every other table already says "not null").
$rng = $d.Content
while ($rng.Find.Execute("id bigint,")) {
    $commaRng = $d.Range($rng.End - 1, $rng.End)
    $commaRng.Text = " not null,"
    $rng.Collapse(0)
}

# --- 2) primary-key constraints: collapse down to (id) ---------------
$d.Content.Find.Execute(" primary key (album_id),", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " primary key (id),", 2) | Out-Null

$d.Content.Find.Execute(" primary key (track_id),", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " primary key (id),", 2) | Out-Null

$d.Content.Find.Execute(" primary key (user_id, album_id)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " primary key (id)", 2) | Out-Null

$d.Content.Find.Execute(" primary key (album_id, user_id)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " primary key (id)", 2) | Out-Null

$d.Content.Find.Execute(" primary key (user_id, track_id)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    " primary key (id)", 2) | Out-Null
